{"js": "// Remove the duplicate \"Ahmad bin Mushtaq\" list entries that follow\n// \"Ali Khalid\" (Electric Power & Propulsion Design Team) and \"Saad bin\n// Tariq\" (Spraying Mechanism), while keeping the \"Ahmad bin Mushtaq\"\n// entry that follows \"Muhammad Adeel\" (UAV Control & Computing System).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst toDelete = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const text = (items[i].text || \"\").trim();\n  if (text !== \"Ahmad bin Mushtaq\") continue;\n\n  const prevText = i > 0 ? (items[i - 1].text || \"\").trim() : \"\";\n  // Keep only the occurrence that directly follows \"Muhammad Adeel\".\n  if (prevText !== \"Muhammad Adeel\") {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the duplicate \"Ahmad bin Mushtaq\" list entries that follow\n# \"Ali Khalid\" (Electric Power & Propulsion Design Team) and \"Saad bin\n# Tariq\" (Spraying Mechanism), while keeping the \"Ahmad bin Mushtaq\"\n# entry that follows \"Muhammad Adeel\" (UAV Control & Computing System).\n\n$d = $word.ActiveDocument\n\n$toDelete = @()\n$prevText = \"\"\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n\n    if ($text -eq \"Ahmad bin Mushtaq\" -and $prevText -ne \"Muhammad Adeel\") {\n        $toDelete += $p\n    }\n\n    $prevText = $text\n}\n\n# Delete from last to first so earlier (still-pending) matches keep\n# pointing at the correct text - deleting forward would shift the\n# document out from under the later saved paragraph references.\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $toDelete[$j].Range.Delete()\n}\n"}
